$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This change represents a new handoff / xliff-generation run: the source
# markdown file's GUID changed (acc6e68f... -> b5178533...), a new xliff
# content hash was produced (abeeb0ae...), the "Latest HO Xliff Generate
# Date" moved from 15:01:27 to 15:01:45, and (because this is a fresh
# handoff) the "Latest Handback" file/date columns on the per-locale sheets
# are no longer populated, so their hyperlink + values are cleared.
# ---------------------------------------------------------------------------

$oldGuid = "acc6e68f-631e-43f6-a6c7-ea97c063043d"
$newGuid = "b5178533-5025-4d14-95f3-e5b78497ead6"
$newHash = "abeeb0ae709d9e9867ee6113bebd454be82a61cf"

$newFileName      = "$newGuid.md"
$newPathAndName   = "e2e\$newGuid.md"
$newGenerateDate  = "2016-08-20 15:01:45"

# Same external hyperlink targets as before (only the visible display text
# and/or underlying value changed in this commit).
$overviewLinkAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/263f9135c1094ede58d26b82153c343b78096741/e2e/$oldGuid.md"
$zhcnLinkAddr      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/263f9135c1094ede58d26b82153c343b78096741/e2e/$oldGuid.md"
$dedeLinkAddr      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/263f9135c1094ede58d26b82153c343b78096741/e2e/$oldGuid.md"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPathAndName
$wsOverview.Range("G2").Value = $newGenerateDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewLinkAddr, "", "", $newPathAndName) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newFileName
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-20 15:01:41"
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

# No handback target/date exists yet for this fresh handoff: clear Latest
# Target File and Latest Handback File (and drop the now-stale hyperlink).
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("I2").Style = "Normal"
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhcnLinkAddr, "", "", $newFileName) | Out-Null

$wsZhCn.Columns.Item(9).ColumnWidth = 17.833333333333336
$wsZhCn.Columns.Item(10).ColumnWidth = 20.833333333333336

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newFileName
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = $newGenerateDate
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("I2").Style = "Normal"
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $dedeLinkAddr, "", "", $newFileName) | Out-Null

$wsDeDe.Columns.Item(9).ColumnWidth = 17.833333333333336
$wsDeDe.Columns.Item(10).ColumnWidth = 20.833333333333336

Write-Output "applied handback report regeneration"
